{"js": "// Collapse the split \"<id>...</id>\" runs (where the inner id text used to\n// sit in its own differently-formatted run) back into a single run per\n// occurrence, e.g. \"<id>\" + \"p082r_a1\" + \"</id>\"  ->  \"<id>p082r_1</id>\".\nconst renames = [\n  [\"p082r_a1\", \"p082r_1\"],\n  [\"p082r_a2\", \"p082r_2\"],\n  [\"p082r_a3\", \"p082r_3\"],\n  [\"p082r_a4\", \"p082r_4\"],\n];\n\nfor (const [oldId, newId] of renames) {\n  const oldText = \"<id>\" + oldId + \"</id>\";\n  const newText = \"<id>\" + newId + \"</id>\";\n\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Collapse the split \"<id>...</id>\" runs (where the inner id text used to\n# sit in its own differently-formatted run) back into a single run per\n# occurrence, e.g. \"<id>\" + \"p082r_a1\" + \"</id>\"  ->  \"<id>p082r_1</id>\".\n$d = $word.ActiveDocument\n\n$renames = @(\n    @{ Old = \"p082r_a1\"; New = \"p082r_1\" },\n    @{ Old = \"p082r_a2\"; New = \"p082r_2\" },\n    @{ Old = \"p082r_a3\"; New = \"p082r_3\" },\n    @{ Old = \"p082r_a4\"; New = \"p082r_4\" }\n)\n\nforeach ($pair in $renames) {\n    $oldText = \"<id>\" + $pair.Old + \"</id>\"\n    $newText = \"<id>\" + $pair.New + \"</id>\"\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}\n"}
